$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'259.88"
$ws.Range("E2").Value = "'1.25%"
$ws.Range("D3").Value = "'26.86"
$ws.Range("E3").Value = "'-1.58%"
$ws.Range("D4").Value = "'4.665"
$ws.Range("E4").Value = "'-0.07%"
$ws.Range("D5").Value = "'0.06063"
$ws.Range("E5").Value = "'3.12%"
$ws.Range("D6").Value = "'6.698"
$ws.Range("E6").Value = "'1.01%"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.8612"
$ws.Range("E7").Value = "'0.35%"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'0.9130"
$ws.Range("E8").Value = "'-2.94%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1403"
$ws.Range("E9").Value = "'0.13%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.05208"
$ws.Range("E10").Value = "'26.12%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07154"
$ws.Range("E11").Value = "'0.71%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03139"
$ws.Range("E12").Value = "'-0.10%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09140"
$ws.Range("E13").Value = "'-0.02%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001536"
$ws.Range("E14").Value = "'0.54%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006065"
$ws.Range("E15").Value = "'0.39%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006120"
$ws.Range("E16").Value = "'-1.61%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.485"
$ws.Range("E17").Value = "'-0.92%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.179"
$ws.Range("E18").Value = "'-0.77%"
$ws.Range("E19").Value = "'-1.27%"
$ws.Range("E20").Value = "'2.44%"
$ws.Range("E21").Value = "'-0.12%"
$ws.Range("E22").Value = "'7.39%"
$ws.Range("D23").Value = "'0.04234"
$ws.Range("E23").Value = "'-0.10%"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("E24").Value = "'-0.60%"
$ws.Range("D25").Value = "'0.004038"
$ws.Range("E25").Value = "'-5.84%"
$ws.Range("E26").Value = "'-0.02%"
$ws.Range("E27").Value = "'-21.37%"
$ws.Range("D40").Value = "'0.03880"
$ws.Range("E40").Value = "'1.46%"
$ws.Range("D41").Value = "'0.1117"
$ws.Range("E41").Value = "'1.30%"
$ws.Range("D42").Value = "'0.004149"
$ws.Range("E42").Value = "'-33.82%"
$ws.Range("E43").Value = "'30.49%"
$ws.Range("D44").Value = "'0.002200"
$ws.Range("E44").Value = "'0.01%"
$ws.Range("E45").Value = "'-2.43%"
$ws.Range("E46").Value = "'0.01%"
$ws.Range("E47").Value = "'9.10%"
$ws.Range("D48").Value = "'0.1353"
$ws.Range("E48").Value = "'-40.60%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.02%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.01%"
